$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.086879134178162
$ws.Range("B1").Value = 2.300782680511475
$ws.Range("C1").Value = 9.712230682373047
$ws.Range("D1").Value = 2.266873121261597
$ws.Range("E1").Value = 1.299554705619812
